$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data rows down by one to make room for a new
# "ECs" sending-cluster row: old row 3 (MuSCs) -> row 4,
# old row 2 (FAPs) -> row 3.
$ws.Range("A3:T3").Copy()
$ws.Range("A4").PasteSpecial(-4104)
$excel.CutCopyMode = 0

$ws.Range("A2:T2").Copy()
$ws.Range("A3").PasteSpecial(-4104)
$excel.CutCopyMode = 0

# --- Row 2: new "ECs" row ---
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cntn2"
$ws.Range("C2").Value = "Cntn1"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.05028033333333334
$ws.Range("H2").Value = 0.150841
$ws.Range("I2").Value = 0.1252715694221136
$ws.Range("J2").Value = 0.1252715694221136
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.293933
$ws.Range("N2").Value = 0.881799
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.01477904921766667
$ws.Range("R2").Value = 0.133011442959
$ws.Range("S2").Value = 0.1252715694221136
$ws.Range("T2").Value = 0.1252715694221136

# --- Row 3: "FAPs" row (previously row 2), with refreshed specificity values ---
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Cntn2"
$ws.Range("C3").Value = "Cntn1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.300673
$ws.Range("H3").Value = 0.9020189999999999
$ws.Range("I3").Value = 0.7491155307811899
$ws.Range("J3").Value = 0.7491155307811898
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.293933
$ws.Range("N3").Value = 0.881799
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.088377716909
$ws.Range("R3").Value = 0.7953994521809999
$ws.Range("S3").Value = 0.7491155307811899
$ws.Range("T3").Value = 0.7491155307811898

# --- Row 4: "MuSCs" row (previously row 3), with refreshed specificity values ---
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Cntn2"
$ws.Range("C4").Value = "Cntn1"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.05041733333333333
$ws.Range("H4").Value = 0.151252
$ws.Range("I4").Value = 0.1256128997966967
$ws.Range("J4").Value = 0.1256128997966967
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.293933
$ws.Range("N4").Value = 0.881799
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.01481931803866667
$ws.Range("R4").Value = 0.133373862348
$ws.Range("S4").Value = 0.1256128997966967
$ws.Range("T4").Value = 0.1256128997966967
